$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new diary entry as row 41
$ws.Range("A41").Value = "4 joulu"
$ws.Range("C41").Value = "Laskentavarjostimen käyttöönottoa"
$ws.Range("B41").Value = "13.00-14.00, 18.00-20.00"
$ws.Range("G41").Value = 3

# Match formatting of the existing rows above (wrap text + time/number styles)
$ws.Range("B41").WrapText = $true
$ws.Range("B41").NumberFormat = "h:mm"
$ws.Range("C41").WrapText = $true

# Row auto-fits to two lines of wrapped text, same as the other entries
$ws.Rows.Item(41).RowHeight = 28.8

# Move the view to the newly-added row, like Excel does after data entry
$ws.Range("H41").Select()
